$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(42613.757210648146, 20, 55, 44, 55, 25, 20739, 17144, 1027, 169, 134, 12, 4, "Named"),
    @(42613.884363425925, 14, 52, 46, 52, 33, 11512, 9142, 559, 90, 80, 10, 5, "Named"),
    @(42614.883136574077, 0, 47, 52, 47, 31, 21435, 21016, 1192, 188, 210, 11, 5, "Named"),
    @(42615.884039351855, 18, 51, 45, 51, 21, 21433, 20343, 1185, 190, 169, 18, 5, "Named")
)

$startRow = 12
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 14; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}
